# Rango proyecciones.xlsx - "probar archivo y arreglar micro errores"
#
# Every data sheet currently has a 5-column header:
#   Material | Pesimista Proy. | Optimista. Proy. | Pesimista Precio | Optimista Precio
# It needs to become a 9-column header:
#   Material | Venta plan | Stock planta | Puerto Chile | Centro Agua |
#   Puerto Oficina | Almacen oficina | Pesimista Proy. | Optimista. Proy.
# with A1:H1 painted bold/white-on-blue, centered + wrapped, thin white border
# (I1, like the old last header cell, keeps the default style).

$wb = $excel.ActiveWorkbook

$headers = @(
    "Material",
    "Venta plan",
    "Stock planta",
    "Puerto Chile",
    "Centro Agua",
    "Puerto Oficina",
    "Almacen oficina",
    "Pesimista Proy.",
    "Optimista. Proy."
)

$BLUE_FILL = 0x00D7A98B   # => fgColor 8BA9D7 in the saved xlsx
$WHITE     = 0x00FFFFFF
$xlCenter  = -4108
$xlPasteFormats = -4122

$styledTemplate = $null

foreach ($ws in $wb.Worksheets) {
    $a1 = $ws.Cells.Item(1, 1).Value2
    if ($a1 -ne "Material") {
        continue
    }

    for ($i = 0; $i -lt $headers.Length; $i++) {
        $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
    }

    $headerRange = $ws.Range("A1:H1")

    if ($null -eq $styledTemplate) {
        # Build the header look exactly once for the whole workbook so the
        # shared style table only gains the font/fill/border/xf combination
        # actually needed - every other header just reuses it via a format
        # copy/paste instead of re-deriving it property by property.
        $headerRange.Font.Bold = $true
        $headerRange.Font.Color = $WHITE
        $headerRange.Interior.Color = $BLUE_FILL
        $headerRange.Borders.LineStyle = 1
        $headerRange.Borders.Color = $WHITE
        $headerRange.HorizontalAlignment = $xlCenter
        $headerRange.VerticalAlignment = $xlCenter
        $headerRange.WrapText = $true

        $styledTemplate = $ws.Range("A1")
    } else {
        $styledTemplate.Copy()
        $headerRange.PasteSpecial($xlPasteFormats)
    }
}

$excel.CutCopyMode = $false
